# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. This updates the DAMSLTag (column I) and
# DialogAct (column J) values for a set of rows in the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 8;   Tag = "ba"; Act = "Appreciation" },
    @{ Row = 10;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 12;  Tag = "ba"; Act = "Appreciation" },
    @{ Row = 19;  Tag = "ba"; Act = "Appreciation" },
    @{ Row = 24;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 53;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 81;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 99;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 108; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 126; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 140; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 205; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 226; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 236; Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 249; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 250; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 256; Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 257; Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 285; Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 308; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 310; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 327; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 338; Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 339; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 348; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 354; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 359; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 364; Tag = "ba"; Act = "Appreciation" },
    @{ Row = 380; Tag = "sd"; Act = "Statement-non-opinion" }
)

foreach ($change in $changes) {
    $ws.Cells.Item($change.Row, 9).Value = $change.Tag
    $ws.Cells.Item($change.Row, 10).Value = $change.Act
}
